$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "29.383.33"); ensure Excel keeps them as
# text instead of auto-converting to numbers / losing trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.383.33"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.874.25"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "0.7117"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "242.06"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.07787"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").Value = "0.3108"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "25.13"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").Value = "0.08450"
$ws.Range("D12").Value = "1.867.72"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "5.238"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "0.7119"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "91.15"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "29.385.05"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "6.043"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000008215"
$ws.Range("E18").Value = "  +5.17%  "
$ws.Range("D19").Value = "241.02"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").Value = "13.26"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "2.121.10"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "7.777"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "163.68"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").Value = "9.050"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").Value = "1.512"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").Value = "4.431"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").Value = "4.302"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "0.05283"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").Value = "1.942"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").Value = "1.178"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "0.7452"
$ws.Range("E36").Value = "  -9.09%  "
$ws.Range("D37").Value = "2.697"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "0.01869"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").Value = "1.215.70"
$ws.Range("E39").Value = "  +5.22%  "
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").Value = "6.484"
$ws.Range("E41").Value = "  +4.31%  "
$ws.Range("D42").Value = "0.8894"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").Value = "72.73"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "108.94"
$ws.Range("E44").Value = "  +6.84%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "2.017.86"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.00000000123"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.5210"
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("D50").Value = "9.368"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "0.4320"
$ws.Range("E51").Value = "  +0.77%  "
